$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 42
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 1
